# ---------------------------------------------------------------------------
# feat: add 2022-Q3 data
#
# 1. Insert a new worksheet "2022-Q3" right before the existing "2022-Q2"
#    sheet (pushing 2022-Q2 / 2022-Q1 / 2021-Q4 / 2021-Q2 one slot to the
#    right) and fill it with the quarterly fund-holdings table.
# 2. Insert a new row at the top of the "总计" (summary) sheet's data table
#    for the 2022-Q3 figures, shifting the older rows down and renumbering
#    the running index column.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: "总计" summary sheet - insert the new 2022-Q3 row (row 2)
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

# Push existing data rows (2..5) down to (3..6), keeping their formatting.
$summary.Range("A2").EntireRow.Insert()

# The freshly-inserted row inherited the header row's style; re-stamp it
# with the regular data-row style (copied from the row right below it).
$summary.Range("A3:D3").Copy()
$summary.Range("A2:D2").PasteSpecial(-4122)

# New 2022-Q3 figures.
$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 11
$summary.Range("D2").Value = 1.31

# Renumber the running index column for the rows that shifted down.
$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3
$summary.Range("A6").Value = 4

# ---------------------------------------------------------------------------
# Step 2: brand-new "2022-Q3" worksheet with the fund holdings table
# ---------------------------------------------------------------------------
$q2Sheet = $wb.Worksheets.Item("2022-Q2")
$newSheet = $wb.Worksheets.Add($q2Sheet)
$newSheet.Name = "2022-Q3"

# --- header row (row 1), formatting copied from the "2022-Q2" header ---
$q2Sheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# --- running-index column (A2:A12), formatting copied from "2022-Q2" ---
$q2Sheet.Range("A2").Copy()
$newSheet.Range("A2:A12").PasteSpecial(-4122)

# Columns B, D, E, F, G hold numeric-looking figures that must stay text
# (e.g. fund codes with leading zeros, "4.40" rather than 4.4). Force the
# text number-format before writing so COM does not coerce them to numbers.
$newSheet.Range("B2:B12").NumberFormat = "@"
$newSheet.Range("D2:D12").NumberFormat = "@"
$newSheet.Range("E2:E12").NumberFormat = "@"
$newSheet.Range("F2:F12").NumberFormat = "@"
$newSheet.Range("G2:G12").NumberFormat = "@"

$fundData = @(
    @("000924", "宝盈先进制造灵活配置混合A",       "7.79", "87.93", "4.40", "0.3428", 10),
    @("013890", "国泰睿毅三年持有期混合A",         "4.86", "89.26", "5.10", "0.2479", 1),
    @("008227", "宝盈研究精选混合A",               "3.98", "89.60", "5.89", "0.2344", 7),
    @("005585", "银河文体娱乐主题灵活配置混合A",   "3.01", "90.28", "6.98", "0.2101", 1),
    @("001628", "招商体育文化休闲股票A",           "2.23", "92.42", "5.06", "0.1128", 5),
    @("008228", "宝盈研究精选混合C",               "0.95", "89.60", "5.89", "0.0560", 7),
    @("015667", "银河文体娱乐主题灵活配置混合C",   "0.41", "90.28", "6.98", "0.0286", 1),
    @("013891", "国泰睿毅三年持有期混合C",         "0.45", "89.26", "5.10", "0.0230", 1),
    @("007579", "宝盈先进制造灵活配置混合C",       "0.52", "87.93", "4.40", "0.0229", 10),
    @("003397", "银华体育文化灵活配置混合",         "0.32", "81.07", "4.47", "0.0143", 6),
    @("015395", "招商体育文化休闲股票C",           "0.25", "92.42", "5.06", "0.0126", 5)
)

$rowIdx = 2
foreach ($row in $fundData) {
    $newSheet.Cells.Item($rowIdx, 1).Value = ($rowIdx - 2)
    $newSheet.Cells.Item($rowIdx, 2).Value = $row[0]
    $newSheet.Cells.Item($rowIdx, 3).Value = $row[1]
    $newSheet.Cells.Item($rowIdx, 4).Value = $row[2]
    $newSheet.Cells.Item($rowIdx, 5).Value = $row[3]
    $newSheet.Cells.Item($rowIdx, 6).Value = $row[4]
    $newSheet.Cells.Item($rowIdx, 7).Value = $row[5]
    $newSheet.Cells.Item($rowIdx, 8).Value = $row[6]
    $rowIdx = $rowIdx + 1
}

# Drop the temporary text-format styling (keeps values as text, but clears
# the now-unneeded explicit number format / quote-prefix styling) so the
# cells match the plain, unstyled data cells used elsewhere in the workbook.
$newSheet.Range("B2:B12").Style = "Normal"
$newSheet.Range("D2:D12").Style = "Normal"
$newSheet.Range("E2:E12").Style = "Normal"
$newSheet.Range("F2:F12").Style = "Normal"
$newSheet.Range("G2:G12").Style = "Normal"

# Match the page margins used throughout the rest of the workbook
# (0.75in / 1in / 0.5in --> 54pt / 72pt / 36pt).
$newSheet.PageSetup.LeftMargin = 54
$newSheet.PageSetup.RightMargin = 54
$newSheet.PageSetup.TopMargin = 72
$newSheet.PageSetup.BottomMargin = 72
$newSheet.PageSetup.HeaderMargin = 36
$newSheet.PageSetup.FooterMargin = 36

# Restore the original active sheet/selection so the only structural change
# to bookViews is the new sheet tab itself.
$summary.Activate()
$summary.Range("A1").Select()

Write-Host "2022-Q3 sheet + summary row added"
